$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.9458723333333334
$ws.Range("H2").Value = 2.837617
$ws.Range("I2").Value = 0.1873686327665471
$ws.Range("J2").Value = 0.1873686327665471
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6124329999999999
$ws.Range("N2").Value = 1.837299
$ws.Range("Q2").Value = 0.5792834307203333
$ws.Range("R2").Value = 5.213550876483
$ws.Range("S2").Value = 0.1873686327665471
$ws.Range("T2").Value = 0.1873686327665471

# Row 3
$ws.Range("I3").Value = 0.3891165466060174
$ws.Range("J3").Value = 0.3891165466060174
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6124329999999999
$ws.Range("N3").Value = 1.837299
$ws.Range("Q3").Value = 1.203022964622
$ws.Range("R3").Value = 10.827206681598
$ws.Range("S3").Value = 0.3891165466060174
$ws.Range("T3").Value = 0.3891165466060174

# Row 4
$ws.Range("G4").Value = 1.499502333333333
$ws.Range("H4").Value = 4.498507
$ws.Range("I4").Value = 0.2970376573303378
$ws.Range("J4").Value = 0.2970376573303378
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6124329999999999
$ws.Range("N4").Value = 1.837299
$ws.Range("Q4").Value = 0.9183447125103331
$ws.Range("R4").Value = 8.265102412592999
$ws.Range("S4").Value = 0.2970376573303378
$ws.Range("T4").Value = 0.2970376573303378

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.6384806666666667
$ws.Range("H5").Value = 1.915442
$ws.Range("I5").Value = 0.1264771632970977
$ws.Range("J5").Value = 0.1264771632970977
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6124329999999999
$ws.Range("N5").Value = 1.837299
$ws.Range("Q5").Value = 0.3910266301286666
$ws.Range("R5").Value = 3.519239671158
$ws.Range("S5").Value = 0.1264771632970977
$ws.Range("T5").Value = 0.1264771632970977
